$wb = $excel.ActiveWorkbook

# --- Update "Metadata" sheet: bump the "Last Updated" timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 01:19 PM"

# --- Update "Stock List" sheet: new top row, everything else shifts down ---
$ws = $wb.Worksheets.Item("Stock List")

# Insert a fresh row at position 2 (pushes all existing data rows down by one)
$ws.Rows.Item(2).Insert()

# The inserted row inherits formatting from the row above it (the bold header);
# strip that so it matches the plain style used by every other data row.
$ws.Range("A2:H2").ClearFormats()

$ws.Range("A2").Value = [char]0x1F4CB
$ws.Range("B2").Value = "CAPTRU-RE1"
$ws.Range("C2").Value = "CAPTRU-RE1"
$ws.Range("D2").Value = 5.67
$ws.Range("E2").Value = -11.9565
$ws.Range("F2").Value = "N/A"
$ws.Range("G2").Value = "N/A"
$ws.Range("H2").Value = 0

# The insert duplicated the last row (old data just shifted down); remove the
# now-redundant trailing row so the sheet keeps its original 76-row extent.
$ws.Rows.Item(77).Delete()
